$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ 2 = -17.55301424166521; 3 = 2.345783445939642; 4 = -17.55301424166521; 5 = -17.55301424166521; 6 = -17.55301424166521; 7 = -17.55301424166521; 8 = -17.55301424166521; 9 = -17.55301424166521; 10 = -17.55301424166521; 11 = -17.55301424166521 }
    3 = @{ 2 = -17.55301424166521; 3 = -17.55301424166521; 4 = -17.55301424166521; 5 = -17.55301424166521; 6 = -17.55301424166521; 7 = -17.55301424166521; 8 = -17.55301424166521; 9 = 2.389077823627419; 10 = -17.55301424166521; 11 = -17.55301424166521 }
    4 = @{ 2 = -17.55301424166521; 3 = 2.046203069798591; 4 = 2.911739884498636; 5 = -17.55301424166521; 6 = 2.497446908407591; 7 = -17.55301424166521; 8 = -17.55301424166521; 9 = -17.55301424166521; 10 = 2.384746053248337; 11 = -17.55301424166521 }
    5 = @{ 2 = -17.55301424166521; 3 = 0.8340817193162681; 4 = -17.55301424166521; 5 = -17.55301424166521; 6 = -17.55301424166521; 7 = 2.026156196201311; 8 = -17.55301424166521; 9 = -17.55301424166521; 10 = -17.55301424166521; 11 = -17.55301424166521 }
    6 = @{ 2 = -17.55301424166521; 3 = -17.55301424166521; 4 = -17.55301424166521; 5 = -17.55301424166521; 6 = -17.55301424166521; 7 = -17.55301424166521; 8 = -17.55301424166521; 9 = -17.55301424166521; 10 = -17.55301424166521; 11 = -17.55301424166521 }
    7 = @{ 2 = 2.954467817389004; 3 = -17.55301424166521; 4 = -17.55301424166521; 5 = -17.55301424166521; 6 = -17.55301424166521; 7 = -17.55301424166521; 8 = -17.55301424166521; 9 = -17.55301424166521; 10 = -17.55301424166521; 11 = -17.55301424166521 }
    8 = @{ 2 = -17.55301424166521; 3 = -17.55301424166521; 4 = -17.55301424166521; 5 = 2.838280750578374; 6 = -17.55301424166521; 7 = -17.55301424166521; 8 = -17.55301424166521; 9 = -17.55301424166521; 10 = -17.55301424166521; 11 = -17.55301424166521 }
    9 = @{ 2 = 3.614527022937744; 3 = -17.55301424166521; 4 = -17.55301424166521; 5 = -17.55301424166521; 6 = -17.55301424166521; 7 = -17.55301424166521; 8 = -17.55301424166521; 9 = -17.55301424166521; 10 = -17.55301424166521; 11 = -17.55301424166521 }
    10 = @{ 2 = -17.55301424166521; 3 = -17.55301424166521; 4 = -17.55301424166521; 5 = -17.55301424166521; 6 = -17.55301424166521; 7 = -17.55301424166521; 8 = -17.55301424166521; 9 = 1.729099098032305; 10 = -17.55301424166521; 11 = 2.139058035936045 }
    11 = @{ 2 = -17.55301424166521; 3 = -17.55301424166521; 4 = -17.55301424166521; 5 = 1.975090706494126; 6 = -17.55301424166521; 7 = 2.282272895653312; 8 = -17.55301424166521; 9 = -17.55301424166521; 10 = -17.55301424166521; 11 = 1.347112566206821 }
    12 = @{ 2 = -17.55301424166521; 3 = -17.55301424166521; 4 = -17.55301424166521; 5 = -17.55301424166521; 6 = -17.55301424166521; 7 = -17.55301424166521; 8 = -17.55301424166521; 9 = -17.55301424166521; 10 = -17.55301424166521; 11 = -17.55301424166521 }
    13 = @{ 2 = -17.55301424166521; 3 = -17.55301424166521; 4 = -17.55301424166521; 5 = 1.59546482312217; 6 = -17.55301424166521; 7 = -17.55301424166521; 8 = -17.55301424166521; 9 = -17.55301424166521; 10 = 2.545553785295027; 11 = 1.536376753537522 }
    14 = @{ 2 = -17.55301424166521; 3 = -17.55301424166521; 4 = 1.539270614297082; 5 = -17.55301424166521; 6 = -17.55301424166521; 7 = -17.55301424166521; 8 = -17.55301424166521; 9 = -17.55301424166521; 10 = -17.55301424166521; 11 = 2.003620907960387 }
    15 = @{ 2 = -17.55301424166521; 3 = -17.55301424166521; 4 = -0.1756257582257256; 5 = -17.55301424166521; 6 = -17.55301424166521; 7 = -17.55301424166521; 8 = -17.55301424166521; 9 = -17.55301424166521; 10 = -17.55301424166521; 11 = -17.55301424166521 }
    16 = @{ 2 = -17.55301424166521; 3 = -17.55301424166521; 4 = -17.55301424166521; 5 = -17.55301424166521; 6 = -17.55301424166521; 7 = -17.55301424166521; 8 = -17.55301424166521; 9 = -17.55301424166521; 10 = 2.441285325175015; 11 = -17.55301424166521 }
    17 = @{ 2 = -17.55301424166521; 3 = 0.4867027147359631; 4 = -0.2094034297329429; 5 = -17.55301424166521; 6 = -17.55301424166521; 7 = -17.55301424166521; 8 = -17.55301424166521; 9 = 0.03907045234749695; 10 = 0.7307297821565886; 11 = -17.55301424166521 }
    18 = @{ 2 = -17.55301424166521; 3 = -17.55301424166521; 4 = -17.55301424166521; 5 = -17.55301424166521; 6 = -17.55301424166521; 7 = -17.55301424166521; 8 = 4.321920967745292; 9 = 0.6877817479499798; 10 = 0.8862363381253955; 11 = -17.55301424166521 }
    19 = @{ 2 = -17.55301424166521; 3 = -17.55301424166521; 4 = 1.692893542681524; 5 = -17.55301424166521; 6 = -17.55301424166521; 7 = -17.55301424166521; 8 = -17.55301424166521; 9 = 1.986835589253494; 10 = -17.55301424166521; 11 = -17.55301424166521 }
    20 = @{ 2 = -17.55301424166521; 3 = 1.890125593612963; 4 = 2.196883003107927; 5 = -17.55301424166521; 6 = 3.843285841862963; 7 = -17.55301424166521; 8 = -17.55301424166521; 9 = 2.276369392988475; 10 = -17.55301424166521; 11 = 2.618344611627957 }
    21 = @{ 2 = -17.55301424166521; 3 = 1.961902030148747; 4 = -17.55301424166521; 5 = 2.559441810156422; 6 = -17.55301424166521; 7 = 3.46757850787054; 8 = -17.55301424166521; 9 = -17.55301424166521; 10 = -17.55301424166521; 11 = -17.55301424166521 }
}

foreach ($r in $data.Keys) {
    foreach ($c in $data[$r].Keys) {
        $ws.Cells.Item($r, $c).Value = $data[$r][$c]
    }
}
